# Apply the target edit:
#  - On "component_list", fill the (previously empty) pos_x / pos_y columns
#    (H2:I86) with "NA" for every component row.
#  - Make "component_list" the active/selected sheet (was "comp_type_dmg_algo"),
#    with the selection left on H85.
#  - Leave "comp_type_dmg_algo" with its prior selection (C6), it simply
#    stops being the active tab.

$wb = $excel.ActiveWorkbook

# --- component_list: stamp pos_x / pos_y with "NA" for rows 2-86 ---
$compList = $wb.Worksheets.Item("component_list")
$compList.Range("H2:I86").Value = "NA"

# --- switch the active tab to component_list, update its selection ---
$compList.Activate()
$compList.Range("H85").Select()
